$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A1: date moves forward one month (serial 45406 -> 45436)
$ws.Range("A1").Value = 45436

# D30 / D31: updated prices
$ws.Range("D30").Value = 1576
$ws.Range("D31").Value = 1794

# Re-touch each merged range (unmerge + remerge) so the merge list is
# rewritten in this exact order, matching the refreshed layout.
$ws.Range("A1:D1").UnMerge()
$ws.Range("A1:D1").Merge()

$ws.Range("B30:C30").UnMerge()
$ws.Range("B30:C30").Merge()

$ws.Range("B29:C29").UnMerge()
$ws.Range("B29:C29").Merge()

$ws.Range("A12:D12").UnMerge()
$ws.Range("A12:D12").Merge()

$ws.Range("B31:C31").UnMerge()
$ws.Range("B31:C31").Merge()

$ws.Range("A11:D11").UnMerge()
$ws.Range("A11:D11").Merge()
